$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.797.53'
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').Value = '1.635.66'
$ws.Range('E3').Value = '  +0.33%  '
$ws.Range('E4').Value = '  -0.80%  '
$ws.Range('D5').Value = '216.72'
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('D6').Value = '0.505'
$ws.Range('E6').Value = '  +1.85%  '
$ws.Range('E7').Value = '  -0.82%  '
$ws.Range('E8').Value = '  +1.93%  '
$ws.Range('D9').Value = '0.0624'
$ws.Range('E9').Value = '  +0.85%  '
$ws.Range('D10').Value = '19.75'
$ws.Range('E10').Value = '  +4.69%  '
$ws.Range('D11').Value = '0.0842'
$ws.Range('E11').Value = '  -0.12%  '
$ws.Range('D12').Value = '1.863.88'
$ws.Range('E12').Value = '  +0.23%  '
$ws.Range('D13').Value = '1.627.61'
$ws.Range('E13').Value = '  -0.26%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.10'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.61%  '
$ws.Range('E15').Value = '  +1.52%  '
$ws.Range('D16').Value = '66.13'
$ws.Range('E16').Value = '  +3.41%  '
$ws.Range('D17').Value = '26.800.16'
$ws.Range('E17').Value = '  +0.43%  '
$ws.Range('D18').Value = '0.0₃0727'
$ws.Range('E18').Value = '  +0.44%  '
$ws.Range('D19').Value = '218.14'
$ws.Range('E19').Value = '  +3.13%  '
$ws.Range('E20').Value = '  -0.78%  '
$ws.Range('E21').Value = '  +1.38%  '
$ws.Range('D22').Value = '6.57'
$ws.Range('E22').Value = '  +6.09%  '
$ws.Range('D23').Value = '2.41'
$ws.Range('E23').Value = '  +2.18%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.10'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.40%  '
$ws.Range('D25').Value = '146.51'
$ws.Range('E25').Value = '  -0.62%  '
$ws.Range('E26').Value = '  -0.64%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.40'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +6.47%  '
$ws.Range('E28').Value = '  +1.65%  '
$ws.Range('D29').Value = '15.76'
$ws.Range('E29').Value = '  +1.47%  '
$ws.Range('E30').Value = '  +1.20%  '
$ws.Range('E31').Value = '  -0.74%  '
$ws.Range('D32').Value = '3.35'
$ws.Range('E32').Value = '  +0.12%  '
$ws.Range('E33').Value = '  +1.71%  '
$ws.Range('D34').Value = '1.54'
$ws.Range('E34').Value = '  +2.12%  '
$ws.Range('E35').Value = '  -0.48%  '
$ws.Range('D36').Value = '1.239.02'
$ws.Range('E36').Value = '  -1.79%  '
$ws.Range('E37').Value = '  +0.87%  '
$ws.Range('D38').Value = '0.533'
$ws.Range('E38').Value = '  +2.22%  '
$ws.Range('D39').Value = '0.825'
$ws.Range('E39').Value = '  +3.39%  '
$ws.Range('E40').Value = '  -0.79%  '
$ws.Range('D41').Value = '0.804'
$ws.Range('E41').Value = '  +0.69%  '
$ws.Range('D43').Value = '1.777.85'
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = '61.39'
$ws.Range('E44').Value = '  +3.38%  '
$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D45').Value = '2.08'
$ws.Range('E45').Value = '  -2.82%  '
$ws.Range('D46').Value = '91.37'
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('E47').Value = '  +0.89%  '
$ws.Range('E48').Value = '  +11.36%  '
$ws.Range('D49').Value = '0.0513'
$ws.Range('E49').Value = '  -0.53%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '0.0965'
$ws.Range('E50').Value = '  +1.21%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '7.55'
$ws.Range('E51').Value = '  +2.38%  '
